# wip - testing schema changes
# Adds Date_column (E) and Datetime_column (F) to the Bin_Master sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -------------------------------------------------------------
$ws.Range("E1").Value = "Date_column"
$ws.Range("F1").Value = "Datetime_column"

# --- Datetime_column (F) string values, in row order 4 then 2 so the ---
# --- shared-string table is built in the same order as the target file --
$ws.Range("F4").Value = "2024-01-01T00:00:00"
$ws.Range("F2").Value = "2024-01-01T16:00:00"

# --- Date_column (E): rows 2, 4 & 5 stay blank but get a date format,  ---
# --- row 3 gets an actual (text) date value -------------------------------
$ws.Range("E2").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("E4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("E5").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2024-01-01"

# --- Column E width --------------------------------------------------------
$ws.Range("E1").ColumnWidth = 14.8

# --- Sheet view / selection -------------------------------------------------
$ws.Range("G12").Select()

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1
